# edit.ps1
# Applies the changes described by the target diff:
#  1. Removes 5 duplicate "learn.chtc.wisc.edu" textboxes that had been
#     accidentally appended to slides 16-20 (the one on slide 21 is kept).
#  2. Updates the exercise numbers referenced in the slide-21 title.
#  3. Updates file-size figures in the three comparison tables found on
#     slides 8, 23 and 25 (10GB -> 20GB, 10MB -> 100MB, etc.)

$p = $ppt.ActivePresentation

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. Delete the stray duplicate "learn.chtc.wisc.edu" textboxes.
# ---------------------------------------------------------------------
$staleTextboxes = @(
    @{ Slide = 16; Name = "Google Shape;302;p18" },
    @{ Slide = 17; Name = "Google Shape;323;p19" },
    @{ Slide = 18; Name = "Google Shape;343;p20" },
    @{ Slide = 19; Name = "Google Shape;364;p21" },
    @{ Slide = 20; Name = "Google Shape;385;p22" }
)

foreach ($entry in $staleTextboxes) {
    $slide = $p.Slides.Item($entry.Slide)
    $shape = Get-ShapeByName $slide $entry.Name
    if ($shape -ne $null) {
        $shape.Delete()
    }
}

# ---------------------------------------------------------------------
# 2. Slide 21 title: "Ex. 4.1-4.2" -> "Ex. 3.1-3.2"
# ---------------------------------------------------------------------
$slide21 = $p.Slides.Item(21)
$title21 = $slide21.Shapes.Title
$title21.TextFrame.TextRange.Text = "At UW-Madison (Ex. 3.1-3.2)"

# ---------------------------------------------------------------------
# 3. Table updates.
# ---------------------------------------------------------------------

# -- Slide 8, table shape "Google Shape;181;p10" --------------------
$slide8 = $p.Slides.Item(8)
$tbl181 = (Get-ShapeByName $slide8 "Google Shape;181;p10").Table
$tbl181.Cell(5, 1).Shape.TextFrame.TextRange.Text = "1GB - 20GB, unique or shared file"
$tbl181.Cell(6, 1).Shape.TextFrame.TextRange.Text = "20 GB " + [char]0x2013 + " TBs, unique or shared"

# -- Slide 23, table shape "Google Shape;421;p25" --------------------
$slide23 = $p.Slides.Item(23)
$tbl421 = (Get-ShapeByName $slide23 "Google Shape;421;p25").Table
$tbl421.Cell(5, 1).Shape.TextFrame.TextRange.Text = "1GB - 20GB, unique or shared file"
$tbl421.Cell(6, 1).Shape.TextFrame.TextRange.Text = "20 GB " + [char]0x2013 + " TBs, unique or shared"

# -- Slide 25, table shape "Google Shape;440;p27" --------------------
$slide25 = $p.Slides.Item(25)
$tbl440 = (Get-ShapeByName $slide25 "Google Shape;440;p27").Table
$tbl440.Cell(2, 3).Shape.TextFrame.TextRange.Text = "100 MB/file (in), 1 GB/file (out); 1 GB/tot (either)"

$cell440_43 = $tbl440.Cell(4, 3).Shape.TextFrame.TextRange
$cell440_43.Runs(1, 1).Text = "20 GB/file"

Write-Output "done"
